$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

# The "input" and "status" column headers on the "data" sheet were
# placeholders that needed fixing - rename them (this updates both the
# header cell and the bound table's column name).
$ws.Range("B1").Value = "~InputName~"
$ws.Range("C1").Value = "~RecordType~"

# Leave the selection where the user ended up after editing the headers.
$ws.Range("C2").Select() | Out-Null
